$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Update month/year values used for the ROI date comparison range.
# FromDate month: December -> May (keep the text/quote-prefix cell style)
$ws.Range("B3").Value = "'May"
# ToDate month: February -> June (keep the text/quote-prefix cell style)
$ws.Range("E3").Value = "'June"

# ToDate year: 2020 -> 2019 (reapply original cell formatting afterwards so the
# existing cell style is retained on the rewritten numeric cell)
$ws.Range("F3").Value = 2019
$ws.Range("C3").Copy()
$ws.Range("F3").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Update the selected cell on the sheet
$ws.Range("K4").Select()
